$d = $word.ActiveDocument
$sel = $word.Selection

# ------------------------------------------------------------------
# The pre-existing "_GoBack" bookmark sits mid-document; once we add
# new trailing content it needs to point at the very end of the doc
# (matching real Word's "last edit location" bookmark behaviour).
# Remove the old one now - we'll recreate it after the new text.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Move to the very end of the document body.
$sel.EndKey(6)

# --- "Training" heading -------------------------------------------------
$sel.TypeParagraph()
$sel.TypeText("")
$sel.TypeParagraph()
$sel.TypeText("Training")

# blank separator paragraph
$sel.TypeParagraph()
$sel.TypeText("")

# body paragraph
$sel.TypeParagraph()
$sel.TypeText("`tOur scientists and developers will benefit from a better understanding of what ParaView has to offer via training. We would like to have a general ParaView training session geared toward scientists and an advanced training session geared towards developers. The developer session would focus more on ")
$sel.TypeText("and ")
$sel.TypeText("understanding of the Paraview and VTK code base.")

# blank separator paragraph
$sel.TypeParagraph()
$sel.TypeText("")

# --- "Help Hours" heading -----------------------------------------------
$sel.TypeParagraph()
$sel.TypeText("Help Hours")

# blank separator paragraph
$sel.TypeParagraph()
$sel.TypeText("")

# body paragraph
$sel.TypeParagraph()
$sel.TypeText("`tOur development roadmap has us moving into some key areas (MPI")
$sel.TypeText(" utilization")
$sel.TypeText(") where ta")
$sel.TypeText("rgeted")
$sel.TypeText(" ")
$sel.TypeText("help would decrease the develop")
$sel.TypeText("ment")
$sel.TypeText(" time. We would request a standard block of 25 hours.")

# ------------------------------------------------------------------
# Recreate "_GoBack" as a zero-length bookmark right after the final
# text we just typed (mirrors Word's behaviour of tracking the last
# edited spot). A temporary placeholder character works around the
# fact that Bookmarks.Add on an already-collapsed range anchors to
# the start of the document instead of the given position.
# ------------------------------------------------------------------
$sel.TypeText("X")
$markRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$delRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$delRange.Delete()
